$d = $word.ActiveDocument

# The final paragraph in the body (an otherwise-empty list item, right
# before the trailing blank list item / sectPr) needs the new sentence
# added as three runs - "hemophilia" is split out on its own (it was
# wrapped in spell-check proofErr markers in the source) while keeping
# identical run formatting (Nunito Sans / color 6C6D74 / white shading)
# throughout, matching the formatting already used on the sibling
# paragraph just above ("Source plasma is plasma that is collected...").

$targetPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Grab a run that already carries the exact rFonts/color/shd combo we
# need (rPr with w:shd at the *run* level, not just the paragraph level)
# and copy/paste it into the target paragraph so the new run(s) inherit
# that formatting faithfully.
$formattedSource = $d.Paragraphs.Item(4).Range.Duplicate
$formattedSource.Copy()

$insertionPoint = $targetPara.Range
$insertionPoint.Collapse(1)
$insertionPoint.Paste()

# Now overwrite the pasted run's text with the actual sentence we want.
$part1 = "The plasma protein therapeutics industry supports volunteerism donation in all of its forms. Source plasma donation and blood donation are critically important activities that contribute to saving lives. Source plasma and recovered plasma are used to produce therapies that treat people with rare, chronic diseases and disorders such as primary immunodeficiency, "
$part2 = "hemophilia"
$part3 = " and a genetic lung disease, as well as in the treatment of trauma, burns and shock."
$fullText = $part1 + $part2 + $part3

$targetPara.Range.Text = $fullText

# Force "hemophilia" onto its own run (matching the source document,
# where it sits between <w:proofErr w:type="spellStart"/> /
# <w:proofErr w:type="spellEnd"/> markers) by nudging a character
# property on just that span and then reverting it - this splits the
# run without leaving any residual formatting difference behind.
$paraStart = $targetPara.Range.Start
$hemoStart = $paraStart + $part1.Length
$hemoEnd = $hemoStart + $part2.Length
$hemoRange = $d.Range($hemoStart, $hemoEnd)
$hemoRange.Bold = 1
$hemoRange.Bold = 0
